$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1 (公告明细): remove hyperlinks, clear old data body ----
$ws1.Hyperlinks.Delete()
$ws1.Range("A2:F9").Clear()

# ---- Sheet1: restyle header row to match sheet2 header style (s=3) ----
$ws2.Range("A1").Copy()
$ws1.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Sheet1: write the new full data body (16 rows) ----
$ws1.Range("A2:F17").Clear()
$ws1.Cells.Item(2,1).Value = 'ST时万'
$ws1.Cells.Item(2,2).Value = '''600241'
$ws1.Cells.Item(2,2).Style = 'Normal'
$ws1.Cells.Item(2,3).Value = 'ST时万:辽宁时代万恒股份有限公司关于撤销其他风险警示暨停牌的公告'
$ws1.Cells.Item(2,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(2,5).Value = '2023-05-29 18:47:58:000'
$ws1.Cells.Item(2,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587340569_1.pdf?1685386082000.pdf'
$ws1.Cells.Item(3,1).Value = 'ST大集'
$ws1.Cells.Item(3,2).Value = '''000564'
$ws1.Cells.Item(3,2).Style = 'Normal'
$ws1.Cells.Item(3,3).Value = 'ST大集:2022年年度股东大会决议公告'
$ws1.Cells.Item(3,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(3,5).Value = '2023-05-26 18:06:00:000'
$ws1.Cells.Item(3,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587293532_1.pdf?1685302108000.pdf'
$ws1.Cells.Item(4,1).Value = 'ST大集'
$ws1.Cells.Item(4,2).Value = '''000564'
$ws1.Cells.Item(4,2).Style = 'Normal'
$ws1.Cells.Item(4,3).Value = 'ST大集:关于供销大集集团股份有限公司2022年年度股东大会的法律意见书'
$ws1.Cells.Item(4,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(4,5).Value = '2023-05-26 18:05:41:000'
$ws1.Cells.Item(4,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587293505_1.pdf?1685124360000.pdf'
$ws1.Cells.Item(5,1).Value = '*ST碳元'
$ws1.Cells.Item(5,2).Value = '''603133'
$ws1.Cells.Item(5,2).Style = 'Normal'
$ws1.Cells.Item(5,3).Value = '*ST碳元:碳元科技股份有限公司关于终止设立控股子公司的公告'
$ws1.Cells.Item(5,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(5,5).Value = '2023-05-29 18:04:41:000'
$ws1.Cells.Item(5,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587339143_1.pdf?1685383506000.pdf'
$ws1.Cells.Item(6,1).Value = '*ST碳元'
$ws1.Cells.Item(6,2).Value = '''603133'
$ws1.Cells.Item(6,2).Style = 'Normal'
$ws1.Cells.Item(6,3).Value = '*ST碳元:碳元科技股份有限公司关于拟投资设立控股子公司的公告'
$ws1.Cells.Item(6,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(6,5).Value = '2023-05-29 17:24:31:000'
$ws1.Cells.Item(6,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587337645_1.pdf?1685381096000.pdf'
$ws1.Cells.Item(7,1).Value = '*ST碳元'
$ws1.Cells.Item(7,2).Value = '''603133'
$ws1.Cells.Item(7,2).Style = 'Normal'
$ws1.Cells.Item(7,3).Value = '*ST碳元:碳元科技股份有限公司第四届董事会第七次会议决议公告'
$ws1.Cells.Item(7,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(7,5).Value = '2023-05-29 17:25:04:000'
$ws1.Cells.Item(7,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587337644_1.pdf?1685381096000.pdf'
$ws1.Cells.Item(8,1).Value = '*ST碳元'
$ws1.Cells.Item(8,2).Value = '''603133'
$ws1.Cells.Item(8,2).Style = 'Normal'
$ws1.Cells.Item(8,3).Value = '*ST碳元:碳元科技股份有限公司关于全资子公司对外投资暨签订招商引资合同书的公告'
$ws1.Cells.Item(8,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(8,5).Value = '2023-05-29 17:25:04:000'
$ws1.Cells.Item(8,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587337643_1.pdf?1685381096000.pdf'
$ws1.Cells.Item(9,1).Value = '*ST碳元'
$ws1.Cells.Item(9,2).Value = '''603133'
$ws1.Cells.Item(9,2).Style = 'Normal'
$ws1.Cells.Item(9,3).Value = '*ST碳元:碳元科技股份有限公司关于召开2023年第三次临时股东大会的通知'
$ws1.Cells.Item(9,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(9,5).Value = '2023-05-29 17:25:04:000'
$ws1.Cells.Item(9,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587337642_1.pdf?1685381096000.pdf'
$ws1.Cells.Item(10,1).Value = '*ST莫高'
$ws1.Cells.Item(10,2).Value = '''600543'
$ws1.Cells.Item(10,2).Style = 'Normal'
$ws1.Cells.Item(10,3).Value = '*ST莫高:莫高股份关于2022年度暨2023年第一季度业绩暨分红说明会召开情况的公告'
$ws1.Cells.Item(10,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(10,5).Value = '2023-05-26 15:33:02:000'
$ws1.Cells.Item(10,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587287498_1.pdf?1685124392000.pdf'
$ws1.Cells.Item(11,1).Value = '佳禾智能'
$ws1.Cells.Item(11,2).Value = '''300793'
$ws1.Cells.Item(11,2).Style = 'Normal'
$ws1.Cells.Item(11,3).Value = '佳禾智能:关于特定股东减持时间过半未减持股份的公告'
$ws1.Cells.Item(11,4).Value = '2023-05-26 00:00:00'
$ws1.Cells.Item(11,5).Value = '2023-05-26 17:40:04:000'
$ws1.Cells.Item(11,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587292362_1.pdf?1685122833000.pdf'
$ws1.Cells.Item(12,1).Value = '孚日转债'
$ws1.Cells.Item(12,2).Value = '''002083'
$ws1.Cells.Item(12,2).Style = 'Normal'
$ws1.Cells.Item(12,3).Value = '孚日股份:关于召开2023年第一次临时股东大会的提示性公告'
$ws1.Cells.Item(12,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(12,5).Value = '2023-05-29 15:40:56:000'
$ws1.Cells.Item(12,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587333839_1.pdf?1685394605000.pdf'
$ws1.Cells.Item(13,1).Value = '绿动转债'
$ws1.Cells.Item(13,2).Value = '''601330'
$ws1.Cells.Item(13,2).Style = 'Normal'
$ws1.Cells.Item(13,3).Value = '绿色动力:绿色动力环保集团股份有限公司2023年度跟踪评级报告'
$ws1.Cells.Item(13,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(13,5).Value = '2023-05-26 17:19:28:000'
$ws1.Cells.Item(13,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587291500_1.pdf?1685121792000.pdf'
$ws1.Cells.Item(14,1).Value = '绿动转债'
$ws1.Cells.Item(14,2).Value = '''601330'
$ws1.Cells.Item(14,2).Style = 'Normal'
$ws1.Cells.Item(14,3).Value = '绿色动力:关于可转换公司债券2023年度跟踪评级结果的公告'
$ws1.Cells.Item(14,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(14,5).Value = '2023-05-26 17:17:09:000'
$ws1.Cells.Item(14,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587291499_1.pdf?1685305345000.pdf'
$ws1.Cells.Item(15,1).Value = 'ST广珠'
$ws1.Cells.Item(15,2).Value = '''600382'
$ws1.Cells.Item(15,2).Style = 'Normal'
$ws1.Cells.Item(15,3).Value = 'ST广珠:广东明珠集团股份有限公司关于控股股东及其一致行动人部分股份质押及部分股份解除质押的公告'
$ws1.Cells.Item(15,4).Value = '2023-05-30 00:00:00'
$ws1.Cells.Item(15,5).Value = '2023-05-29 18:49:44:000'
$ws1.Cells.Item(15,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305291587340878_1.pdf?1685386208000.pdf'
$ws1.Cells.Item(16,1).Value = '伊利股份'
$ws1.Cells.Item(16,2).Value = '''600887'
$ws1.Cells.Item(16,2).Style = 'Normal'
$ws1.Cells.Item(16,3).Value = '伊利股份:内蒙古伊利实业集团股份有限公司关于调整2022年度利润分配现金分红总额的公告'
$ws1.Cells.Item(16,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(16,5).Value = '2023-05-26 18:04:49:000'
$ws1.Cells.Item(16,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587293447_1.pdf?1685124532000.pdf'
$ws1.Cells.Item(17,1).Value = '伊利股份'
$ws1.Cells.Item(17,2).Value = '''600887'
$ws1.Cells.Item(17,2).Style = 'Normal'
$ws1.Cells.Item(17,3).Value = '伊利股份:内蒙古伊利实业集团股份有限公司关于股份回购实施结果暨股份变动公告'
$ws1.Cells.Item(17,4).Value = '2023-05-27 00:00:00'
$ws1.Cells.Item(17,5).Value = '2023-05-26 18:04:49:000'
$ws1.Cells.Item(17,6).Value = 'https://pdf.dfcfw.com/pdf/H2_AN202305261587293446_1.pdf?1685124301000.pdf'

# ---- Sheet1: page margins (match sheet2's defaults) ----
$ws1.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws1.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws1.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws1.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

$ws1.Range("A1").Select()

# ---- Sheet2 (公告汇总): clear old data body, write new summary ----
$ws2.Range("A2:C6").Clear()
$ws2.Range("A2:C10").Clear()
$ws2.Cells.Item(2,1).Value = 'ST时万'
$ws2.Cells.Item(2,2).Value = '''600241'
$ws2.Cells.Item(2,2).Style = 'Normal'
$ws2.Cells.Item(2,3).Value = 1
$ws2.Cells.Item(3,1).Value = 'ST大集'
$ws2.Cells.Item(3,2).Value = '''000564'
$ws2.Cells.Item(3,2).Style = 'Normal'
$ws2.Cells.Item(3,3).Value = 2
$ws2.Cells.Item(4,1).Value = '*ST碳元'
$ws2.Cells.Item(4,2).Value = '''603133'
$ws2.Cells.Item(4,2).Style = 'Normal'
$ws2.Cells.Item(4,3).Value = 5
$ws2.Cells.Item(5,1).Value = '*ST莫高'
$ws2.Cells.Item(5,2).Value = '''600543'
$ws2.Cells.Item(5,2).Style = 'Normal'
$ws2.Cells.Item(5,3).Value = 1
$ws2.Cells.Item(6,1).Value = '佳禾智能'
$ws2.Cells.Item(6,2).Value = '''300793'
$ws2.Cells.Item(6,2).Style = 'Normal'
$ws2.Cells.Item(6,3).Value = 1
$ws2.Cells.Item(7,1).Value = '孚日转债'
$ws2.Cells.Item(7,2).Value = '''002083'
$ws2.Cells.Item(7,2).Style = 'Normal'
$ws2.Cells.Item(7,3).Value = 1
$ws2.Cells.Item(8,1).Value = '绿动转债'
$ws2.Cells.Item(8,2).Value = '''601330'
$ws2.Cells.Item(8,2).Style = 'Normal'
$ws2.Cells.Item(8,3).Value = 2
$ws2.Cells.Item(9,1).Value = 'ST广珠'
$ws2.Cells.Item(9,2).Value = '''600382'
$ws2.Cells.Item(9,2).Style = 'Normal'
$ws2.Cells.Item(9,3).Value = 1
$ws2.Cells.Item(10,1).Value = '伊利股份'
$ws2.Cells.Item(10,2).Value = '''600887'
$ws2.Cells.Item(10,2).Style = 'Normal'
$ws2.Cells.Item(10,3).Value = 2

$ws2.Range("A1").Select()

